# Apply cryptos list update (cell text/value changes per commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin name / link swap for rows 28-29 (PEPE <-> WrappedeETH)
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("B29").Value = "PEPE"
$ws.Range("C29").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"

# Price column updates that Excel would not misinterpret as numbers
$ws.Range("D2").Value = "62.430.38"
$ws.Range("D3").Value = "2.460.10"
$ws.Range("D9").Value = "2.458.00"
$ws.Range("D16").Value = "2.913.33"
$ws.Range("D17").Value = "62.346.00"
$ws.Range("D18").Value = "2.460.94"
$ws.Range("D28").Value = "2.584.98"
$ws.Range("D29").Value = "0.0₃0966"
$ws.Range("D51").Value = "0.0₆0239"

# Price column updates that look like plain numbers -- force text storage
# (NumberFormat "@" while assigning, then ClearFormats so the cell keeps no
# explicit style, same as the original inline-string cells) so the value is
# preserved exactly as text instead of being parsed into a Double.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D4").Value = "1.00"
$ws.Range("D5").Value = "583.65"
$ws.Range("D6").Value = "144.14"
$ws.Range("D13").Value = "0.346"
$ws.Range("D14").Value = "26.67"
$ws.Range("D19").Value = "10.96"
$ws.Range("D20").Value = "7.19"
$ws.Range("D21").Value = "331.54"
$ws.Range("D23").Value = "2.02"
$ws.Range("D25").Value = "66.12"
$ws.Range("D26").Value = "9.46"
$ws.Range("D27").Value = "629.38"
$ws.Range("D32").Value = "8.06"
$ws.Range("D35").Value = "4.96"
$ws.Range("D39").Value = "5.37"
$ws.Range("D40").Value = "149.91"
$ws.Range("D41").Value = "18.41"
$ws.Range("D42").Value = "1.75"
$ws.Range("D43").Value = "42.51"
$ws.Range("D46").Value = "144.45"
$ws.Range("D48").Value = "0.0529"
$ws.Range("D49").Value = "0.604"
$ws.Range("D50").Value = "19.76"

$ws.Range("D4").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D21").ClearFormats()
$ws.Range("D23").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D35").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D41").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D50").ClearFormats()

# Volume(1h) percentage text updates
$ws.Range("E2").Value = "  -1.45%  "
$ws.Range("E3").Value = "  +0.16%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E6").Value = "  -1.83%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  -3.54%  "
$ws.Range("E11").Value = "  +2.06%  "
$ws.Range("E12").Value = "  -1.02%  "
$ws.Range("E13").Value = "  -2.93%  "
$ws.Range("E14").Value = "  -1.62%  "
$ws.Range("E15").Value = "  -4.09%  "
$ws.Range("E16").Value = "  +0.37%  "
$ws.Range("E17").Value = "  -1.41%  "
$ws.Range("E18").Value = "  +0.61%  "
$ws.Range("E19").Value = "  -3.20%  "
$ws.Range("E20").Value = "  -1.89%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  -2.14%  "
$ws.Range("E23").Value = "  -3.77%  "
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("E26").Value = "  +5.11%  "
$ws.Range("E27").Value = "  +1.47%  "
$ws.Range("E28").Value = "  +0.82%  "
$ws.Range("E29").Value = "  -6.39%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("E31").Value = "  -4.53%  "
$ws.Range("E32").Value = "  -2.66%  "
$ws.Range("E33").Value = "  +0.77%  "
$ws.Range("E34").Value = "  -0.97%  "
$ws.Range("E35").Value = "  -4.94%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E37").Value = "  -6.19%  "
$ws.Range("E38").Value = "  -0.36%  "
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  +2.03%  "
$ws.Range("E41").Value = "  -2.54%  "
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("E45").Value = "  -4.88%  "
$ws.Range("E46").Value = "  -2.98%  "
$ws.Range("E47").Value = "  -3.00%  "
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("E50").Value = "  -7.07%  "
$ws.Range("E51").Value = "  +8.51%  "
